# Updates cryptos list data (prices / 1h volume %) per the Wed Jul 17
# 08:59:35 UTC 2024 GitHub Actions refresh, including the Stacks/Maker
# (rows 37-38) and Cosmos/Stellar (rows 49-50) ranking swaps.
#
# Price-column (D) values that look like plain decimals (e.g. "579.13")
# are written with a leading "'" (quote-prefix) so Excel keeps them as
# Text, matching the source workbook's inline-string cells, instead of
# silently converting them to numbers (which would also mangle values
# like "1.00" -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.211.71'
$ws.Range("E2").Value = '  +3.82%  '

$ws.Range("D3").Value = '3.482.50'
$ws.Range("E3").Value = '  +3.24%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '''579.13'
$ws.Range("E5").Value = '  +3.09%  '

$ws.Range("D6").Value = '''162.04'
$ws.Range("E6").Value = '  +5.04%  '

$ws.Range("D7").Value = '''0.615'
$ws.Range("E7").Value = '  +14.01%  '

$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '3.483.65'
$ws.Range("E9").Value = '  +3.28%  '

$ws.Range("E10").Value = '  -1.45%  '

$ws.Range("D11").Value = '''0.125'
$ws.Range("E11").Value = '  +4.23%  '

$ws.Range("D12").Value = '''0.446'
$ws.Range("E12").Value = '  +3.68%  '

$ws.Range("D13").Value = '4.084.19'
$ws.Range("E13").Value = '  +3.29%  '

$ws.Range("E14").Value = '  +0.79%  '

$ws.Range("E15").Value = '  +3.25%  '

$ws.Range("D16").Value = '''28.72'
$ws.Range("E16").Value = '  +6.72%  '

$ws.Range("D17").Value = '65.219.76'
$ws.Range("E17").Value = '  +3.73%  '

$ws.Range("D18").Value = '3.524.98'
$ws.Range("E18").Value = '  +7.95%  '

$ws.Range("D19").Value = '''6.44'
$ws.Range("E19").Value = '  +3.56%  '

$ws.Range("E20").Value = '  +2.80%  '

$ws.Range("D21").Value = '''383.37'
$ws.Range("E21").Value = '  +2.33%  '

$ws.Range("E22").Value = '  +3.18%  '

$ws.Range("E23").Value = '  +4.95%  '

$ws.Range("E24").Value = '  +2.59%  '

$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("E26").Value = '  +5.44%  '

$ws.Range("D27").Value = '''10.11'
$ws.Range("E27").Value = '  +7.23%  '

$ws.Range("D28").Value = '''0.179'
$ws.Range("E28").Value = '  +1.94%  '

$ws.Range("E29").Value = '  -0.87%  '

$ws.Range("D30").Value = '''1.52'
$ws.Range("E30").Value = '  +13.63%  '

$ws.Range("D31").Value = '''6.25'
$ws.Range("E31").Value = '  +3.15%  '

$ws.Range("E32").Value = '  +4.01%  '

$ws.Range("D33").Value = '''23.67'
$ws.Range("E33").Value = '  +2.87%  '

$ws.Range("D34").Value = '''7.21'
$ws.Range("E34").Value = '  +7.04%  '

$ws.Range("E35").Value = '  +10.43%  '

$ws.Range("D36").Value = '''161.84'
$ws.Range("E36").Value = '  +1.78%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '3.052.45'
$ws.Range("E37").Value = '  +3.59%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '''1.92'
$ws.Range("E38").Value = '  +6.24%  '

$ws.Range("D39").Value = '''0.0774'
$ws.Range("E39").Value = '  +1.44%  '

$ws.Range("D40").Value = '''26.98'
$ws.Range("E40").Value = '  -0.22%  '

$ws.Range("E41").Value = '  +6.84%  '

$ws.Range("D42").Value = '''0.0321'
$ws.Range("E42").Value = '  +1.56%  '

$ws.Range("D43").Value = '''6.57'
$ws.Range("E43").Value = '  +0.49%  '

$ws.Range("D44").Value = '''42.84'
$ws.Range("E44").Value = '  +3.53%  '

$ws.Range("D45").Value = '''0.779'
$ws.Range("E45").Value = '  +5.24%  '

$ws.Range("D46").Value = '''25.91'
$ws.Range("E46").Value = '  +12.60%  '

$ws.Range("D47").Value = '''1.11'
$ws.Range("E47").Value = '  +5.40%  '

$ws.Range("D48").Value = '''317.25'
$ws.Range("E48").Value = '  +11.64%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '''0.111'
$ws.Range("E49").Value = '  +8.15%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''6.76'
$ws.Range("E50").Value = '  +6.60%  '

$ws.Range("D51").Value = '''2.21'
$ws.Range("E51").Value = '  +5.85%  '
